$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(45, 8).Value = 1100
$ws.Cells.Item(45, 9).Value = 1000
$ws.Cells.Item(45, 10).Value = 1200
$ws.Cells.Item(45, 11).Value = 3000
$ws.Cells.Item(45, 12).Value = 3600
$ws.Cells.Item(45, 13).Value = -2808
$ws.Cells.Item(45, 14).Value = -3984
$ws.Cells.Item(98, 8).Value = 825
$ws.Cells.Item(98, 10).Value = 0
$ws.Cells.Item(98, 12).Value = 0
$ws.Cells.Item(98, 14).Value = ""
$ws.Cells.Item(107, 8).Value = 666.62067
$ws.Cells.Item(107, 9).Value = 563.7857
$ws.Cells.Item(107, 10).Value = 3546
$ws.Cells.Item(107, 11).Value = 563.7857
$ws.Cells.Item(107, 12).Value = 3546
$ws.Cells.Item(107, 13).Value = 1356.2143
$ws.Cells.Item(107, 14).Value = -7386
$ws.Cells.Item(122, 8).Value = 825
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 14).Value = ""
$ws.Cells.Item(125, 8).Value = 616
$ws.Cells.Item(125, 9).Value = 616
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 11).Value = 5544
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 13).Value = ""
$ws.Cells.Item(125, 14).Value = -3084
$ws.Cells.Item(128, 8).Value = 10000
$ws.Cells.Item(128, 10).Value = 10000
$ws.Cells.Item(128, 12).Value = 10000
$ws.Cells.Item(128, 14).Value = -19960
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9531.487999999999
$ws.Cells.Item(32, 9).Value = 6653.684
$ws.Cells.Item(32, 10).Value = 31402.8
$ws.Cells.Item(32, 11).Value = 6653.684
$ws.Cells.Item(32, 12).Value = 31402.8
$ws.Cells.Item(32, 13).Value = -6366.684
$ws.Cells.Item(32, 14).Value = -31976.8
$ws.Cells.Item(61, 8).Value = 2887.92
$ws.Cells.Item(61, 9).Value = 3126.0952
$ws.Cells.Item(61, 11).Value = 3126.0952
$ws.Cells.Item(61, 13).Value = -2914.0952
$ws.Cells.Item(132, 8).Value = 24044.348
$ws.Cells.Item(132, 9).Value = 2151.1
$ws.Cells.Item(132, 10).Value = 169999.33
$ws.Cells.Item(132, 11).Value = 6453.299999999999
$ws.Cells.Item(132, 12).Value = 509997.99
$ws.Cells.Item(132, 13).Value = -3923.299999999999
$ws.Cells.Item(132, 14).Value = -515057.99
$ws.Cells.Item(136, 8).Value = 2887.92
$ws.Cells.Item(136, 9).Value = 3126.0952
$ws.Cells.Item(136, 11).Value = 9378.285600000001
$ws.Cells.Item(136, 13).Value = -6828.285600000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 3266.4595
$ws.Cells.Item(134, 9).Value = 3612.8708
$ws.Cells.Item(134, 11).Value = 10838.6124
$ws.Cells.Item(134, 13).Value = -8303.6124
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3495.0588
$ws.Cells.Item(31, 10).Value = 4663.0625
$ws.Cells.Item(31, 12).Value = 4663.0625
$ws.Cells.Item(31, 14).Value = -5253.0625
$ws.Cells.Item(34, 8).Value = 3495.0588
$ws.Cells.Item(34, 10).Value = 4663.0625
$ws.Cells.Item(34, 12).Value = 4663.0625
$ws.Cells.Item(34, 14).Value = -5067.0625
$ws.Cells.Item(58, 8).Value = 16665.727
$ws.Cells.Item(58, 9).Value = 1344.1428
$ws.Cells.Item(58, 10).Value = 43478.5
$ws.Cells.Item(58, 11).Value = 1344.1428
$ws.Cells.Item(58, 12).Value = 43478.5
$ws.Cells.Item(58, 13).Value = -1141.1428
$ws.Cells.Item(58, 14).Value = -43884.5
$ws.Cells.Item(99, 8).Value = 23813068
$ws.Cells.Item(99, 10).Value = 41670336
$ws.Cells.Item(99, 12).Value = 41670336
$ws.Cells.Item(99, 14).Value = -41673332
$ws.Cells.Item(122, 8).Value = 1438.1333
$ws.Cells.Item(122, 9).Value = 1357.2
$ws.Cells.Item(122, 10).Value = 1600
$ws.Cells.Item(122, 11).Value = 4071.6
$ws.Cells.Item(122, 12).Value = 4800
$ws.Cells.Item(122, 13).Value = -1621.6
$ws.Cells.Item(122, 14).Value = -9700
$ws.Cells.Item(126, 8).Value = 23813068
$ws.Cells.Item(126, 10).Value = 41670336
$ws.Cells.Item(126, 12).Value = 125011008
$ws.Cells.Item(126, 14).Value = -125015948
$ws.Cells.Item(127, 8).Value = 37446.5
$ws.Cells.Item(127, 10).Value = 37446.5
$ws.Cells.Item(127, 12).Value = 37446.5
$ws.Cells.Item(127, 14).Value = -47366.5
$ws.Cells.Item(132, 8).Value = 5301.3335
$ws.Cells.Item(132, 9).Value = 4272
$ws.Cells.Item(132, 10).Value = 6742.4
$ws.Cells.Item(132, 11).Value = 12816
$ws.Cells.Item(132, 12).Value = 20227.2
$ws.Cells.Item(132, 13).Value = -10286
$ws.Cells.Item(132, 14).Value = -25287.2
$ws.Cells.Item(134, 8).Value = 1303.92
$ws.Cells.Item(134, 9).Value = 1221.2858
$ws.Cells.Item(134, 10).Value = 1409.091
$ws.Cells.Item(134, 11).Value = 3663.8574
$ws.Cells.Item(134, 12).Value = 4227.272999999999
$ws.Cells.Item(134, 13).Value = -1128.8574
$ws.Cells.Item(134, 14).Value = -9297.272999999999
$ws.Cells.Item(136, 8).Value = 16665.727
$ws.Cells.Item(136, 9).Value = 1344.1428
$ws.Cells.Item(136, 10).Value = 43478.5
$ws.Cells.Item(136, 11).Value = 4032.4284
$ws.Cells.Item(136, 12).Value = 130435.5
$ws.Cells.Item(136, 13).Value = -1482.4284
$ws.Cells.Item(136, 14).Value = -135535.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 84.25
$ws.Cells.Item(2, 9).Value = 34.666668
$ws.Cells.Item(2, 10).Value = 114
$ws.Cells.Item(2, 11).Value = 208.000008
$ws.Cells.Item(2, 12).Value = 684
$ws.Cells.Item(2, 13).Value = -95.00000800000001
$ws.Cells.Item(2, 14).Value = -910
$ws.Cells.Item(45, 8).Value = 1000
$ws.Cells.Item(45, 9).Value = 0
$ws.Cells.Item(45, 10).Value = 1000
$ws.Cells.Item(45, 11).Value = 0
$ws.Cells.Item(45, 12).Value = ""
$ws.Cells.Item(45, 13).Value = 3000
$ws.Cells.Item(45, 14).Value = -4064
$ws.Cells.Item(92, 8).Value = 560.7692
$ws.Cells.Item(92, 9).Value = 250
$ws.Cells.Item(92, 10).Value = 1058
$ws.Cells.Item(92, 11).Value = 750
$ws.Cells.Item(92, 12).Value = 3174
$ws.Cells.Item(92, 13).Value = 498
$ws.Cells.Item(92, 14).Value = -5670
$ws.Cells.Item(97, 8).Value = 367
$ws.Cells.Item(97, 9).Value = 222.5
$ws.Cells.Item(97, 10).Value = 463.33334
$ws.Cells.Item(97, 11).Value = 667.5
$ws.Cells.Item(97, 12).Value = 1390.00002
$ws.Cells.Item(97, 13).Value = -171.5
$ws.Cells.Item(97, 14).Value = -2382.00002
$ws.Cells.Item(129, 8).Value = 233007.31
$ws.Cells.Item(129, 10).Value = 465384.62
$ws.Cells.Item(129, 12).Value = 1396153.86
$ws.Cells.Item(129, 14).Value = -1406153.86
$ws.Cells.Item(131, 8).Value = 106045.75
$ws.Cells.Item(131, 9).Value = 698
$ws.Cells.Item(131, 10).Value = 111898.4
$ws.Cells.Item(131, 11).Value = 2094
$ws.Cells.Item(131, 12).Value = 335695.2
$ws.Cells.Item(131, 13).Value = 2946
$ws.Cells.Item(131, 14).Value = -345775.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(24, 8).Value = 200000
$ws.Cells.Item(24, 9).Value = 200000
$ws.Cells.Item(24, 10).Value = 200000
$ws.Cells.Item(24, 11).Value = 200000
$ws.Cells.Item(24, 12).Value = 200000
$ws.Cells.Item(24, 13).Value = -199827
$ws.Cells.Item(24, 14).Value = -200346
$ws.Cells.Item(126, 8).Value = 5515.357
$ws.Cells.Item(126, 9).Value = 4487.5
$ws.Cells.Item(126, 10).Value = 6885.8335
$ws.Cells.Item(126, 11).Value = 13462.5
$ws.Cells.Item(126, 12).Value = 20657.5005
$ws.Cells.Item(126, 13).Value = -10992.5
$ws.Cells.Item(126, 14).Value = -25597.5005
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 1772.4242
$ws.Cells.Item(132, 9).Value = 1288.9412
$ws.Cells.Item(132, 11).Value = 3866.8236
$ws.Cells.Item(132, 13).Value = -1336.8236
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 1458.5714
$ws.Cells.Item(6, 9).Value = 1005
$ws.Cells.Item(6, 10).Value = 1640
$ws.Cells.Item(6, 11).Value = 1005
$ws.Cells.Item(6, 12).Value = 1640
$ws.Cells.Item(6, 13).Value = -890
$ws.Cells.Item(6, 14).Value = -1870
$ws.Cells.Item(70, 8).Value = 20000
$ws.Cells.Item(70, 10).Value = 20000
$ws.Cells.Item(70, 12).Value = 20000
$ws.Cells.Item(70, 14).Value = -20630
$ws.Cells.Item(73, 8).Value = 20000
$ws.Cells.Item(73, 10).Value = 20000
$ws.Cells.Item(73, 12).Value = 20000
$ws.Cells.Item(73, 14).Value = -22184
$ws.Cells.Item(138, 8).Value = 42714.5
$ws.Cells.Item(138, 10).Value = 50429
$ws.Cells.Item(138, 12).Value = 50429
$ws.Cells.Item(138, 14).Value = -60709
